$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-12 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-13 Tuesday", 2)

$d.Content.Find.Execute("499÷5=99, 4", $true, $false, $false, $false, $false, $true, 1, $false, "337÷3=112, 1", 2)
$d.Content.Find.Execute("585÷2=292, 1", $true, $false, $false, $false, $false, $true, 1, $false, "806÷6=134, 2", 2)
$d.Content.Find.Execute("340÷4=85, 0", $true, $false, $false, $false, $false, $true, 1, $false, "573÷7=81, 6", 2)
$d.Content.Find.Execute("944÷6=157, 2", $true, $false, $false, $false, $false, $true, 1, $false, "782÷5=156, 2", 2)
$d.Content.Find.Execute("429÷8=53, 5", $true, $false, $false, $false, $false, $true, 1, $false, "529÷3=176, 1", 2)
$d.Content.Find.Execute("915÷8=114, 3", $true, $false, $false, $false, $false, $true, 1, $false, "659÷4=164, 3", 2)
$d.Content.Find.Execute("807÷9=89, 6", $true, $false, $false, $false, $false, $true, 1, $false, "681÷8=85, 1", 2)
$d.Content.Find.Execute("462÷9=51, 3", $true, $false, $false, $false, $false, $true, 1, $false, "842÷6=140, 2", 2)
$d.Content.Find.Execute("871÷4=217, 3", $true, $false, $false, $false, $false, $true, 1, $false, "971÷2=485, 1", 2)
$d.Content.Find.Execute("729÷3=243, 0", $true, $false, $false, $false, $false, $true, 1, $false, "335÷8=41, 7", 2)
$d.Content.Find.Execute("299÷3=99, 2", $true, $false, $false, $false, $false, $true, 1, $false, "289÷6=48, 1", 2)
$d.Content.Find.Execute("125÷3=41, 2", $true, $false, $false, $false, $false, $true, 1, $false, "788÷8=98, 4", 2)
$d.Content.Find.Execute("767÷7=109, 4", $true, $false, $false, $false, $false, $true, 1, $false, "828÷7=118, 2", 2)
$d.Content.Find.Execute("778÷4=194, 2", $true, $false, $false, $false, $false, $true, 1, $false, "221÷8=27, 5", 2)
$d.Content.Find.Execute("492÷4=123, 0", $true, $false, $false, $false, $false, $true, 1, $false, "406÷4=101, 2", 2)
$d.Content.Find.Execute("695÷7=99, 2", $true, $false, $false, $false, $false, $true, 1, $false, "954÷7=136, 2", 2)
$d.Content.Find.Execute("538÷2=269, 0", $true, $false, $false, $false, $false, $true, 1, $false, "285÷5=57, 0", 2)
$d.Content.Find.Execute("571÷3=190, 1", $true, $false, $false, $false, $false, $true, 1, $false, "568÷3=189, 1", 2)
$d.Content.Find.Execute("188÷7=26, 6", $true, $false, $false, $false, $false, $true, 1, $false, "502÷4=125, 2", 2)
$d.Content.Find.Execute("793÷9=88, 1", $true, $false, $false, $false, $false, $true, 1, $false, "870÷7=124, 2", 2)
$d.Content.Find.Execute("264÷3=88, 0", $true, $false, $false, $false, $false, $true, 1, $false, "623÷9=69, 2", 2)
$d.Content.Find.Execute("239÷7=34, 1", $true, $false, $false, $false, $false, $true, 1, $false, "405÷9=45, 0", 2)
$d.Content.Find.Execute("417÷4=104, 1", $true, $false, $false, $false, $false, $true, 1, $false, "234÷8=29, 2", 2)
$d.Content.Find.Execute("290÷5=58, 0", $true, $false, $false, $false, $false, $true, 1, $false, "834÷2=417, 0", 2)
$d.Content.Find.Execute("433÷6=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "289÷6=48, 1", 2)
